$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 724.75
$ws.Range("I32").Value = 600
$ws.Range("J32").Value = 766.3333
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 766.3333
$ws.Range("M32").Value = -274
$ws.Range("N32").Value = -1418.3333

$ws.Range("H62").Value = 7608.8335
$ws.Range("I62").Value = 3130
$ws.Range("J62").Value = 30003
$ws.Range("K62").Value = 3130
$ws.Range("L62").Value = 30003
$ws.Range("M62").Value = -2506
$ws.Range("N62").Value = -31251

$ws.Range("H65").Value = 7608.8335
$ws.Range("I65").Value = 3130
$ws.Range("J65").Value = 30003
$ws.Range("K65").Value = 15650
$ws.Range("L65").Value = 150015
$ws.Range("M65").Value = -12530
$ws.Range("N65").Value = -156255

$ws.Range("H132").Value = 1465.8889
$ws.Range("I132").Value = 1387.6538
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 4162.9614
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -1632.9614
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12326.648
$ws.Range("I61").Value = 9380.825999999999
$ws.Range("J61").Value = 17166.215
$ws.Range("K61").Value = 9380.825999999999
$ws.Range("L61").Value = 17166.215
$ws.Range("M61").Value = -9168.825999999999
$ws.Range("N61").Value = -17590.215

$ws.Range("H74").Value = 6995.7144
$ws.Range("I74").Value = 2226.3572
$ws.Range("K74").Value = 2226.3572
$ws.Range("M74").Value = -1352.3572

$ws.Range("H77").Value = 6995.7144
$ws.Range("I77").Value = 2226.3572
$ws.Range("K77").Value = 11131.786
$ws.Range("M77").Value = -6763.786

$ws.Range("H110").Value = 2042.6
$ws.Range("J110").Value = 2142
$ws.Range("L110").Value = 2142
$ws.Range("N110").Value = -6232

$ws.Range("H121").Value = 49999.715
$ws.Range("J121").Value = 49999.715
$ws.Range("L121").Value = 49999.715
$ws.Range("N121").Value = -53493.715

$ws.Range("H132").Value = 3697.1765
$ws.Range("I132").Value = 3753
$ws.Range("K132").Value = 11259
$ws.Range("M132").Value = -8729

$ws.Range("H136").Value = 12326.648
$ws.Range("I136").Value = 9380.825999999999
$ws.Range("J136").Value = 17166.215
$ws.Range("K136").Value = 28142.478
$ws.Range("L136").Value = 51498.645
$ws.Range("M136").Value = -25592.478
$ws.Range("N136").Value = -56598.645

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 10045
$ws.Range("I44").Value = 10045
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 10045
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -9548
$ws.Range("N44").ClearContents()

$ws.Range("H86").Value = 1988.102
$ws.Range("I86").Value = 1796.1063
$ws.Range("J86").Value = 6500
$ws.Range("K86").Value = 1796.1063
$ws.Range("L86").Value = 6500
$ws.Range("M86").Value = -673.1062999999999
$ws.Range("N86").Value = -8746

$ws.Range("H89").Value = 1988.102
$ws.Range("I89").Value = 1796.1063
$ws.Range("J89").Value = 6500
$ws.Range("K89").Value = 8980.531499999999
$ws.Range("L89").Value = 32500
$ws.Range("M89").Value = -3364.531499999999
$ws.Range("N89").Value = -43732

$ws.Range("H134").Value = 37857.535
$ws.Range("I134").Value = 2222.111
$ws.Range("J134").Value = 1000014
$ws.Range("K134").Value = 6666.333
$ws.Range("L134").Value = 3000042
$ws.Range("M134").Value = -4131.333
$ws.Range("N134").Value = -3005112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7080.7676
$ws.Range("I31").Value = 7629.5483
$ws.Range("J31").Value = 5663.0835
$ws.Range("K31").Value = 7629.5483
$ws.Range("L31").Value = 5663.0835
$ws.Range("M31").Value = -7334.5483
$ws.Range("N31").Value = -6253.0835

$ws.Range("H34").Value = 7080.7676
$ws.Range("I34").Value = 7629.5483
$ws.Range("J34").Value = 5663.0835
$ws.Range("K34").Value = 7629.5483
$ws.Range("L34").Value = 5663.0835
$ws.Range("M34").Value = -7427.5483
$ws.Range("N34").Value = -6067.0835

$ws.Range("H58").Value = 2530546.2
$ws.Range("I58").Value = 3954590
$ws.Range("J58").Value = 11084.462
$ws.Range("K58").Value = 3954590
$ws.Range("L58").Value = 11084.462
$ws.Range("M58").Value = -3954387
$ws.Range("N58").Value = -11490.462

$ws.Range("H94").Value = 1003
$ws.Range("I94").Value = 597
$ws.Range("K94").Value = 597
$ws.Range("M94").Value = -146

$ws.Range("H122").Value = 16327.454
$ws.Range("I122").Value = 7287.5557
$ws.Range("K122").Value = 21862.6671
$ws.Range("M122").Value = -19412.6671

$ws.Range("H136").Value = 2530546.2
$ws.Range("I136").Value = 3954590
$ws.Range("J136").Value = 11084.462
$ws.Range("K136").Value = 11863770
$ws.Range("L136").Value = 33253.386
$ws.Range("M136").Value = -11861220
$ws.Range("N136").Value = -38353.386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 588
$ws.Range("I18").Value = 405.6
$ws.Range("J18").Value = 1500
$ws.Range("K18").Value = 1216.8
$ws.Range("L18").Value = 4500
$ws.Range("M18").Value = -1047.8
$ws.Range("N18").Value = -4838

$ws.Range("H38").Value = 62.75
$ws.Range("I38").Value = 33
$ws.Range("J38").Value = 92.5
$ws.Range("K38").Value = 99
$ws.Range("L38").Value = 277.5
$ws.Range("M38").Value = 248
$ws.Range("N38").Value = -971.5

$ws.Range("H40").Value = 170.4
$ws.Range("I40").Value = 88
$ws.Range("J40").Value = 500
$ws.Range("K40").Value = 352
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -283
$ws.Range("N40").Value = -2138

$ws.Range("H60").Value = 343.2
$ws.Range("I60").Value = 194
$ws.Range("K60").Value = 582
$ws.Range("M60").Value = -331

$ws.Range("H92").Value = 500
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -252
$ws.Range("N92").ClearContents()

$ws.Range("H110").Value = 3950
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 3950
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 11850
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -20030

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 12073.583
$ws.Range("I132").Value = 9029.888999999999
$ws.Range("K132").Value = 27089.667
$ws.Range("M132").Value = -24559.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 911.5
$ws.Range("I16").Value = 893.8
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 893.8
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -723.8
$ws.Range("N16").Value = -1340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H81").Value = 14288519
$ws.Range("I81").Value = 1199.5
$ws.Range("J81").Value = 20003448
$ws.Range("K81").Value = 2399
$ws.Range("L81").Value = 40006896
$ws.Range("M81").Value = -1338
$ws.Range("N81").Value = -40009018

$ws.Range("H84").Value = 14288519
$ws.Range("I84").Value = 1199.5
$ws.Range("J84").Value = 20003448
$ws.Range("K84").Value = 11995
$ws.Range("L84").Value = 200034480
$ws.Range("M84").Value = -6691
$ws.Range("N84").Value = -200045088

$ws.Range("H100").Value = 1585.7858
$ws.Range("I100").Value = 569.8
$ws.Range("J100").Value = 4125.75
$ws.Range("K100").Value = 1139.6
$ws.Range("L100").Value = 8251.5
$ws.Range("M100").Value = -598.5999999999999
$ws.Range("N100").Value = -9333.5

$ws.Range("H113").Value = 713.0357
$ws.Range("I113").Value = 297.23077
$ws.Range("J113").Value = 1073.4
$ws.Range("K113").Value = 891.69231
$ws.Range("L113").Value = 3220.2
$ws.Range("M113").Value = 1278.30769
$ws.Range("N113").Value = -7560.200000000001

$ws.Range("H120").Value = 69700
$ws.Range("J120").Value = 69700
$ws.Range("L120").Value = 69700
$ws.Range("N120").Value = -79376

$ws.Range("H136").Value = 4381.721
$ws.Range("I136").Value = 1696.1621
$ws.Range("J136").Value = 8521.958000000001
$ws.Range("K136").Value = 5088.4863
$ws.Range("L136").Value = 25565.874
$ws.Range("M136").Value = -2538.4863
$ws.Range("N136").Value = -30665.874
